$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 10; $r -le 17; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "Data Nutrition Label") {
        $cell.Value = "Dataset Nutrition Label"
    }
}
